# Add the "eluaat" concept row to the materiaalklasse sheet.
#
# The source list is alphabetically sorted by notation; "eluaat" sorts
# between "eieren" (row 14) and "fruit" (row 15). We insert a new row at
# 15 (shifting every following concept row down by one, and updating the
# sheet dimension automatically), fill it with the eluaat concept data,
# and patch the handful of cells elsewhere that reference the set of
# materiaalklasse members/concepts so "eluaat" is included in them too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15, pushing fruit..materiaalklasse (old rows 15-39)
# down to rows 16-40.
$ws.Rows(15).Insert()

# Populate the new row 15 with the "eluaat" concept.
$ws.Cells.Item(15, 1).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/eluaat"
$ws.Cells.Item(15, 2).Value = "http://www.w3.org/2004/02/skos/core#Concept"
$ws.Cells.Item(15, 3).Value = "be.vlaanderen.bodemenondergrond.data.id.concept.materiaalklasse.eluaat"
$ws.Cells.Item(15, 4).Value = "https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/materiaalklasse"
$ws.Cells.Item(15, 5).Value = "null"
$ws.Cells.Item(15, 6).Value = "eluaat"
$ws.Cells.Item(15, 7).Value = "Eluaat"
$ws.Cells.Item(15, 8).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/vastdeelvdaarde"
$ws.Cells.Item(15, 9).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/vastdeelvdaarde"
$ws.Cells.Item(15, 10).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/vastdeelvdaarde"
$ws.Cells.Item(15, 11).Value = "https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/materiaalklasse"
$ws.Cells.Item(15, 12).Value = "null"
$ws.Cells.Item(15, 13).Value = "null"
$ws.Cells.Item(15, 14).Value = "null"
$ws.Cells.Item(15, 15).Value = "null"
$ws.Cells.Item(15, 16).Value = "null"

# Row 2: the "materiaalklasses" collection's member list gains eluaat.
$ws.Cells.Item(2, 5).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/afvalwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/antropogeneobjecten|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bemalingswater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bodem|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bodemlucht|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bodemvocht|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/depositie|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/dierlijkmateriaal|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/dnapl|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/drinkwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/effluent|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/eieren|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/eluaat|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/fruit|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/gftafval|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/groenten|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/grondwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/hardgesteente|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/koelwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/lnapl|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/lucht|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/melk|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/migratie|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/oppervlaktewater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/plantaardigmateriaal|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/puurproduct|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/regenwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/schelpdieren|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/sediment|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/strooisel|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/vastdeelvdaarde|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/vilt|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/water|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/waterbodem|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/zeeschuim|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/zeewater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/zwevendestof"

# Row 33 (formerly row 32, "vastdeelvdaarde") is eluaat's new broader
# concept, so its semanticRelation/narrower/narrowerTransitive lists
# (columns J/L/M) gain eluaat too.
$vastdeelList = "https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bodem|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/eluaat|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/hardgesteente|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/sediment|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/waterbodem"
$ws.Cells.Item(33, 10).Value = $vastdeelList
$ws.Cells.Item(33, 12).Value = $vastdeelList
$ws.Cells.Item(33, 13).Value = $vastdeelList

# Row 40 (formerly row 39, the conceptscheme itself) has its
# hasTopConcept list (column O) gain eluaat too.
$ws.Cells.Item(40, 15).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/afvalwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/antropogeneobjecten|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bemalingswater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bodem|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bodemlucht|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/bodemvocht|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/depositie|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/dierlijkmateriaal|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/dnapl|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/drinkwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/effluent|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/eieren|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/eluaat|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/fruit|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/gftafval|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/groenten|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/grondwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/hardgesteente|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/koelwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/lnapl|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/lucht|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/melk|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/migratie|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/oppervlaktewater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/plantaardigmateriaal|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/puurproduct|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/regenwater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/schelpdieren|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/sediment|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/strooisel|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/vastdeelvdaarde|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/vilt|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/water|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/waterbodem|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/zeeschuim|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/zeewater|https://data.bodemenondergrond.vlaanderen.be/id/concept/materiaalklasse/zwevendestof"

"eluaat row inserted"
